$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Starting SoC (%) value swap
$ws.Range("B6").Value = 99

# Row 7: Ending SoC (%) value swap
$ws.Range("B7").Value = 9

# Row 8: label update
$ws.Range("A8").Value = "Total distance covered (km)"

# Row 9: label update
$ws.Range("A9").Value = "Total energy consumption(WH/KM)"

# Row 10: label + value update
$ws.Range("A10").Value = "Total SOC consumed(%)"
$ws.Range("B10").Value = 90

# Row 12: label update
$ws.Range("A12").Value = "Peak Power(kW)"

# Row 13: label update
$ws.Range("A13").Value = "Average Power(kW)"

# Row 14: label update
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"

# Row 15: label + value update
$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 5.786160427063244

# Row 16: label + value update (now Highest Cell Voltage)
$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.379

# Row 17: label + value update (now Lowest Cell Voltage)
$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 2.93

# Row 18: label update
$ws.Range("A18").Value = "Difference in Cell Voltage(V)"

# Row 19: label update
$ws.Range("A19").Value = "Minimum Temperature(C)"

# Row 20: label update
$ws.Range("A20").Value = "Maximum Temperature(C)"

# Row 21: label + value update
$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 18

# Row 22: label update
$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"

# Row 23: label update
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"

# Row 24: label update
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"

# Row 25: label update
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"

# Row 26: label update
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"

# Row 27: label update
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"

# Row 28: label update (now highest cell temp)
$ws.Range("A28").Value = "highest cell temp(C)"

# Row 29: label update (now lowest cell temp)
$ws.Range("A29").Value = "lowest cell temp(C)"

# Row 30: label update
$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

# Row 31: label + value update (now Battery Voltage)
$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 55

# Row 32: label + value update (now Total energy charged)
$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.908872808333333

# Row 33: label + value update (now Electricity consumption units)
$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = [double]"1.277384839217681e-07"

# Row 34: label + value update (now Idling time percentage)
$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 6.664304335720975

# Row 35: label + value update (now Time spent in 0-10 km/h)
$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 10.19008327216584

# Row 36: label + value update (now Time spent in 10-20 km/h)
$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 5.044420258675238

# Row 37: label + value update (now Time spent in 20-30 km/h)
$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 9.165000379660331

# Row 38: label + value update (now Time spent in 30-40 km/h)
$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 23.28836467640285

# Row 39: label + value update (now Time spent in 40-50 km/h)
$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 18.38315320559872

# Row 40: label + value update (now Time spent in 50-60 km/h)
$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 17.40616062163052

# Row 41: label + value update (now Time spent in 60-70 km/h)
$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 9.127034346604571

# Row 42: label + value update (now Time spent in 70-80 km/h)
$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 0.3999088815206662

# Row 43: new row
$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
